$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to Text format before assigning, then reset the
    # cell style to Normal/General so the stored format matches the
    # original (unstyled) cells while keeping the exact literal text
    # (avoids Excel auto-converting numeric-looking strings, which would
    # drop trailing zeros or switch to scientific notation).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("B2") "Bitcoin"
Set-TextValue $ws.Range("C2") "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextValue $ws.Range("D2") "23.333.61"
Set-TextValue $ws.Range("E2") "  +0.47%  "

Set-TextValue $ws.Range("B3") "Ethereum"
Set-TextValue $ws.Range("C3") "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextValue $ws.Range("D3") "1.628.23"
Set-TextValue $ws.Range("E3") "  +1.23%  "

Set-TextValue $ws.Range("B4") "TetherUSD"
Set-TextValue $ws.Range("C4") "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  +0.38%  "

Set-TextValue $ws.Range("B5") "USDC"
Set-TextValue $ws.Range("C5") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D5") "1.004"
Set-TextValue $ws.Range("E5") "  +0.47%  "

Set-TextValue $ws.Range("B6") "BNB"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue $ws.Range("D6") "303.13"
Set-TextValue $ws.Range("E6") "  -0.76%  "

Set-TextValue $ws.Range("B7") "XRP"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D7") "0.3759"
Set-TextValue $ws.Range("E7") "  +0.05%  "

Set-TextValue $ws.Range("B8") "OKB"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D8") "52.29"
Set-TextValue $ws.Range("E8") "  -1.56%  "

Set-TextValue $ws.Range("B9") "Cardano"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.3608"
Set-TextValue $ws.Range("E9") "  -0.23%  "

Set-TextValue $ws.Range("B10") "Polygon"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D10") "1.228"
Set-TextValue $ws.Range("E10") "  -3.04%  "

Set-TextValue $ws.Range("B11") "BinanceUSD"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D11") "1.006"
Set-TextValue $ws.Range("E11") "  +0.60%  "

Set-TextValue $ws.Range("B12") "Dogecoin"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D12") "0.08055"
Set-TextValue $ws.Range("E12") "  -1.13%  "

Set-TextValue $ws.Range("B13") "Solana"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D13") "22.56"
Set-TextValue $ws.Range("E13") "  -1.35%  "

Set-TextValue $ws.Range("B14") "Polkadot"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "6.550"
Set-TextValue $ws.Range("E14") "  -0.85%  "

Set-TextValue $ws.Range("B15") "ShibaInu"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D15") "0.00001247"
Set-TextValue $ws.Range("E15") "  +0.05%  "

Set-TextValue $ws.Range("B16") "Chainlink"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "7.203"
Set-TextValue $ws.Range("E16") "  -2.20%  "

Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "1.633.07"
Set-TextValue $ws.Range("E17") "  +1.63%  "

Set-TextValue $ws.Range("B18") "Litecoin"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D18") "93.58"
Set-TextValue $ws.Range("E18") "  -0.65%  "

Set-TextValue $ws.Range("B19") "TRON"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D19") "0.06930"
Set-TextValue $ws.Range("E19") "  +0.11%  "

Set-TextValue $ws.Range("B20") "Avalanche"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D20") "17.90"
Set-TextValue $ws.Range("E20") "  -1.72%  "

Set-TextValue $ws.Range("B21") "Dai"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D21") "1.004"
Set-TextValue $ws.Range("E21") "  +0.23%  "

Set-TextValue $ws.Range("B22") "Uniswap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "6.460"
Set-TextValue $ws.Range("E22") "  -1.30%  "

Set-TextValue $ws.Range("B23") "WrappedBTC"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D23") "23.321.34"
Set-TextValue $ws.Range("E23") "  +0.46%  "

Set-TextValue $ws.Range("B24") "Cosmos"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D24") "12.70"
Set-TextValue $ws.Range("E24") "  -1.49%  "

Set-TextValue $ws.Range("B25") "LidoDAOToken"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D25") "3.203"
Set-TextValue $ws.Range("E25") "  +3.47%  "

Set-TextValue $ws.Range("B26") "Toncoin"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D26") "2.425"
Set-TextValue $ws.Range("E26") "  +0.36%  "

Set-TextValue $ws.Range("B27") "EthereumClassic"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "21.00"
Set-TextValue $ws.Range("E27") "  -1.00%  "

Set-TextValue $ws.Range("B28") "Monero"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "148.11"
Set-TextValue $ws.Range("E28") "  -1.70%  "

Set-TextValue $ws.Range("B29") "HuobiToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D29") "5.289"
Set-TextValue $ws.Range("E29") "  +0.28%  "

Set-TextValue $ws.Range("B30") "BitcoinCash"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D30") "134.60"
Set-TextValue $ws.Range("E30") "  -0.45%  "

Set-TextValue $ws.Range("B31") "WEMIXTOKEN"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D31") "2.305"
Set-TextValue $ws.Range("E31") "  -3.66%  "

Set-TextValue $ws.Range("B32") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D32") "1.812.55"
Set-TextValue $ws.Range("E32") "  +1.41%  "

Set-TextValue $ws.Range("B33") "Filecoin"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "6.752"
Set-TextValue $ws.Range("E33") "  -0.62%  "

Set-TextValue $ws.Range("B34") "FraxShare"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D34") "10.83"
Set-TextValue $ws.Range("E34") "  +4.73%  "

Set-TextValue $ws.Range("B35") "ImmutableX"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D35") "0.9469"
Set-TextValue $ws.Range("E35") "  -1.05%  "

Set-TextValue $ws.Range("B36") "VeChain"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D36") "0.02823"
Set-TextValue $ws.Range("E36") "  +1.95%  "

Set-TextValue $ws.Range("B37") "Algorand"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D37") "0.2524"
Set-TextValue $ws.Range("E37") "  +0.37%  "

Set-TextValue $ws.Range("B38") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D38") "6.116"
Set-TextValue $ws.Range("E38") "  -0.34%  "

Set-TextValue $ws.Range("B39") "Stellar"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D39") "0.08806"
Set-TextValue $ws.Range("E39") "  +0.32%  "

Set-TextValue $ws.Range("B40") "Hedera"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.07147"
Set-TextValue $ws.Range("E40") "  -3.03%  "

Set-TextValue $ws.Range("B41") "TrustWalletToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D41") "1.362"
Set-TextValue $ws.Range("E41") "  -2.35%  "

Set-TextValue $ws.Range("B42") "TheSandbox"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.7015"
Set-TextValue $ws.Range("E42") "  -1.15%  "

Set-TextValue $ws.Range("B43") "EnergySwap"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D43") "16.15"
Set-TextValue $ws.Range("E43") "  +1.89%  "

Set-TextValue $ws.Range("B44") "Aptos"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D44") "12.27"
Set-TextValue $ws.Range("E44") "  -1.70%  "

Set-TextValue $ws.Range("B45") "Decentraland"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.6442"
Set-TextValue $ws.Range("E45") "  -1.41%  "

Set-TextValue $ws.Range("B46") "Frax"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D46") "1.003"
Set-TextValue $ws.Range("E46") "  +0.45%  "

Set-TextValue $ws.Range("B47") "NEARProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D47") "2.310"
Set-TextValue $ws.Range("E47") "  -0.95%  "

Set-TextValue $ws.Range("B48") "PancakeSwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D48") "3.980"
Set-TextValue $ws.Range("E48") "  -0.85%  "

Set-TextValue $ws.Range("B49") "Cronos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.07973"
Set-TextValue $ws.Range("E49") "  +0.07%  "

Set-TextValue $ws.Range("B50") "Flow"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
Set-TextValue $ws.Range("D50") "1.205"
Set-TextValue $ws.Range("E50") "  +0.45%  "

Set-TextValue $ws.Range("B51") "Quant"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D51") "126.08"
Set-TextValue $ws.Range("E51") "  -5.64%  "
